$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell E71: "Pendiente ADM" -> "ICD30461848" ---
$ws.Range("E71").Value = "ICD30461848"

# --- Helper: write a value as TEXT (preserve leading digits / date-looking
#     strings exactly as typed) without leaving a lingering number format
#     on the cell, then restore the plain "Normal" style. ---
function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- New row 73 ---
Set-TextCell "A73" "5467"
Set-TextCell "B73" "8/22/2025"
Set-TextCell "C73" "FRANCO 2515"
Set-TextCell "D73" "12"
Set-TextCell "E73" "809102560"
Set-TextCell "F73" "NEW"
Set-TextCell "G73" "Pendiente"
Set-TextCell "H73" "PIcada"
$ws.Range("I73").Value = 1
Set-TextCell "J73" "Cambio"
Set-TextCell "K73" "Sin equipos"
Set-TextCell "L73" "Pasante"
$ws.Range("M73").Value = -58.502342
$ws.Range("N73").Value = -34.578839
Set-TextCell "O73" "Paternal"
Set-TextCell "P73" "Capital Norte"

# --- New row 74 ---
Set-TextCell "A74" "5503"
Set-TextCell "B74" "8/22/2025"
Set-TextCell "C74" "CRAMER AV. 1771"
Set-TextCell "D74" "13"
Set-TextCell "E74" "809102564"
Set-TextCell "F74" "NEW"
Set-TextCell "G74" "Pendiente"
Set-TextCell "H74" "PIcada"
$ws.Range("I74").Value = 1
Set-TextCell "J74" "Cambio"
Set-TextCell "K74" "Sin equipos"
Set-TextCell "L74" "Pasante"
$ws.Range("M74").Value = -58.458506
$ws.Range("N74").Value = -34.56783
Set-TextCell "O74" "Colegiales"
Set-TextCell "P74" "Capital Norte"
